# Append 5 new daily rows (234-238) to Sheet1, continuing the series
# that currently ends at row 233 (date serial 44307 / 2021-04-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44308, 1, 7, 99.5732574679943),
    @(44309, 2, 7, 99.5732574679943),
    @(44310, 0, 7, 99.5732574679943),
    @(44311, 1, 5, 71.12375533428165),
    @(44312, 1, 5, 71.12375533428165)
)

$lastRow = 233
$startRow = 234

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Clone the formatting of the last existing data row (borders, bold,
    # center/top alignment and the custom date/time number format on
    # column A; plain/general formatting on B, C, D) onto the new row.
    $ws.Range("A$lastRow`:D$lastRow").Copy() | Out-Null
    $ws.Range("A$r`:D$r").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    $lastRow = $r
}
